$d = $word.ActiveDocument

# --- locate the run that needs to be split -------------------------------
# Original single run:
#   "Il sistema deve fornire un meccanismo di acquisto di un prodotto
#    attraverso selezionando un apposito pulsante che si troverà nei
#    dettagli dell'annuncio "
$oldText = "Il sistema deve fornire un meccanismo di acquisto di un prodotto attraverso selezionando un apposito pulsante che si troverà nei dettagli dell’annuncio "

$target = $d.Content
$found = $target.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target run text"
}

$oldStart = $target.Start
$oldEnd = $target.End
$oldLen = $oldEnd - $oldStart

# --- build the replacement runs (same rPr: sz=22 / szCs=22) --------------
# "Il sistema deve fornire " + "il" + " meccanismo " + "per acquistare" +
# " un prodotto selezionando un apposito pulsante che si troverà nei
#  dettagli dell'annuncio" + "."
$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Il sistema deve fornire </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>il</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> meccanismo </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>per acquistare</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> un prodotto selezionando un apposito pulsante che si troverà nei dettagli dell’annuncio</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

# Compute the length (in characters) of the text the new runs contain, so
# we can find the (shifted) old run afterwards without relying on a Range
# object automatically re-basing itself across the insertion.
$newParts = @(
    "Il sistema deve fornire ",
    "il",
    " meccanismo ",
    "per acquistare",
    " un prodotto selezionando un apposito pulsante che si troverà nei dettagli dell’annuncio",
    "."
)
$newLen = 0
foreach ($p in $newParts) { $newLen += $p.Length }

# --- insert the new runs just before the old run, then drop the old one --
$insertionPoint = $d.Range($oldStart, $oldStart)
$insertionPoint.Collapse(1)
$insertionPoint.InsertXML($newXml)

$shiftedOldStart = $oldStart + $newLen
$shiftedOldEnd = $oldEnd + $newLen
$staleRun = $d.Range($shiftedOldStart, $shiftedOldEnd)
if ($staleRun.Text -ne $oldText) {
    throw "Stale-run sanity check failed before delete"
}
$staleRun.Delete()

Write-Output "Split run into $($newParts.Count) runs (old length $oldLen -> new length $newLen)"
